$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing values in row 2 (M2:T2) ---
$ws.Range("M2").Value = 0.3730536666666667
$ws.Range("N2").Value = 1.119161
$ws.Range("O2").Value = 0.2333846330299144
$ws.Range("P2").Value = 0.2333846330299144
$ws.Range("Q2").Value = 0.04104137478711111
$ws.Range("R2").Value = 0.369372373084
$ws.Range("S2").Value = 0.2333846330299144
$ws.Range("T2").Value = 0.2333846330299144

# --- Update existing values in row 3 (O3, P3, S3, T3) ---
$ws.Range("O3").Value = 0.4556015723565537
$ws.Range("P3").Value = 0.4556015723565538
$ws.Range("S3").Value = 0.4556015723565537
$ws.Range("T3").Value = 0.4556015723565538

# --- Update existing values in row 4 (M4:T4) ---
$ws.Range("M4").Value = 0.461934
$ws.Range("N4").Value = 1.385802
$ws.Range("O4").Value = 0.2889887078106916
$ws.Range("P4").Value = 0.2889887078106916
$ws.Range("Q4").Value = 0.050819515032
$ws.Range("R4").Value = 0.457375635288
$ws.Range("S4").Value = 0.2889887078106916
$ws.Range("T4").Value = 0.2889887078106916

# --- Add new row 5 ---
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cort"
$ws.Range("C5").Value = "Sstr2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1100146666666667
$ws.Range("H5").Value = 0.330044
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.03520599999999999
$ws.Range("N5").Value = 0.105618
$ws.Range("O5").Value = 0.02202508680284024
$ws.Range("P5").Value = 0.02202508680284025
$ws.Range("Q5").Value = 0.003873176354666666
$ws.Range("R5").Value = 0.034858587192
$ws.Range("S5").Value = 0.02202508680284024
$ws.Range("T5").Value = 0.02202508680284025
